$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values for columns E and F across rows 1-3 (header + two data rows)
$e1 = $ws.Range("E1").Value2
$e2 = $ws.Range("E2").Value2
$e3 = $ws.Range("E3").Value2

$f1 = $ws.Range("F1").Value2
$f2 = $ws.Range("F2").Value2
$f3 = $ws.Range("F3").Value2

# Swap columns E and F
$ws.Range("E1").Value = $f1
$ws.Range("E2").Value = $f2
$ws.Range("E3").Value = $f3

$ws.Range("F1").Value = $e1
$ws.Range("F2").Value = $e2
$ws.Range("F3").Value = $e3

# Update the active selection
$ws.Range("H7").Select()
